$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Have a carefully [thought out] UI that is easy to navigate"
#    Collapse the three runs (with gramStart/gramEnd proofErr markers around
#    "thought out") into a single plain run with the full sentence.
#    A Find/Replace over the whole visible text re-merges the runs and drops
#    the now-redundant proofErr markers, while leaving the paragraph itself
#    (and its identity attributes) untouched.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Have a carefully thought out UI that is easy to navigate", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Have a carefully thought out UI that is easy to navigate", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop the stray <w:lastRenderedPageBreak/> in front of
#    "Split the program into the following sections:" - a no-op text
#    Find/Replace causes the run to be rebuilt without that leftover marker.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Split the program into the following sections:", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Split the program into the following sections:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Using Git and Github" -> give "Git" its own spell-check run, matching
#    the existing treatment of "Github":
#      Using  |<spellStart>Git<spellEnd>|  and  |<spellStart>Github<spellEnd>
#    Plain text Find/Replace cannot introduce new <w:proofErr/> markers, and
#    editing a sub-range of the paragraph leaves the original proofErr
#    elements orphaned at the wrong spot. Rebuilding the whole paragraph via
#    Range.InsertXML (while preserving its own identity attributes, read
#    back from WordOpenXML) is reliable and leaves everything else alone.
# ---------------------------------------------------------------------------
function Get-ParagraphOpenTagAttrs($needleText) {
    $full = $d.Content.WordOpenXML
    $pos = $full.IndexOf($needleText)
    $prefix = $full.Substring(0, $pos)
    $pStart = $prefix.LastIndexOf("<w:p ")
    $openEnd = $full.IndexOf(">", $pStart)
    $openTag = $full.Substring($pStart, $openEnd + 1 - $pStart)
    $m = [regex]::Match($openTag, '^<w:p ([^>]*)>')
    return $m.Groups[1].Value
}

$gitAttrs = Get-ParagraphOpenTagAttrs("Using Git and Github")

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Using Git and Github") {
        $gitXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' +
            '<w:p ' + $gitAttrs + '>' +
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr>' +
            '<w:r><w:t xml:space="preserve">Using </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '</w:p>' +
            '</w:body></w:document>' +
            '</pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($gitXml) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 4) Add a new bullet after "... temporarily installed?)" (still before the
#    trailing blank paragraph at the end of the document) describing backing
#    up the writable partition with rsync.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*temporarily installed?)*") {
        $p.Range.InsertParagraphAfter() | Out-Null
    }
}

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*temporarily installed?)*") {
        $newPara = $p.Next()
        $newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' +
            '<w:p>' +
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="4"/></w:numPr></w:pPr>' +
            '<w:r><w:t xml:space="preserve">It may be best to have the deploy script back up the current contents of </w:t></w:r>' +
            '<w:r><w:t>the writable partition to a folder in the project on the laptop; probably using rsync</w:t></w:r>' +
            '</w:p>' +
            '</w:body></w:document>' +
            '</pkg:xmlData></pkg:part></pkg:package>'
        $newPara.Range.InsertXML($newXml) | Out-Null
    }
}
